$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: column C label changes from "speed" to "fps" ---
# (A1/B1/D1 text stays the same; setting C1 causes the now-unused
# "speed" shared string to drop out and the table to re-pack.)
$ws.Range("C1").Value = "fps"

# --- Grade-table data (A2:D31) ---
# Each row: distanceFurlong(A), speed/fps(C), grade(D)
$data = @(
    @(6, 56.8, 6),
    @(6, 55.8, 5),
    @(6, 54.8, 4),
    @(6, 53.8, 3),
    @(6, 52.8, 2),
    @(5, 57, 6),
    @(5, 56, 5),
    @(5, 55, 4),
    @(5, 54, 3),
    @(5, 53, 2),
    @(4, 57.2, 6),
    @(4, 56.2, 5),
    @(4, 55.2, 4),
    @(4, 54.2, 3),
    @(4, 53.2, 2),
    @(3, 57.5, 6),
    @(3, 56.5, 5),
    @(3, 55.5, 4),
    @(3, 54.5, 3),
    @(3, 53.5, 2),
    @(2, 58.5, 6),
    @(2, 57.5, 5),
    @(2, 56.5, 4),
    @(2, 55.5, 3),
    @(2, 54.5, 2),
    @(1, 59, 6),
    @(1, 58, 5),
    @(1, 57, 4),
    @(1, 56, 3),
    @(1, 55, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}

# --- Selection moves to E10 ---
$ws.Range("E10").Select()
